$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Column D (Price) updates
Set-TextValue $ws 'D2' '30.184.06'
Set-TextValue $ws 'D3' '1.850.02'
Set-TextValue $ws 'D4' '1.002'
Set-TextValue $ws 'D5' '235.41'
Set-TextValue $ws 'D7' '0.4701'
Set-TextValue $ws 'D8' '0.2888'
Set-TextValue $ws 'D9' '0.06521'
Set-TextValue $ws 'D10' '21.70'
Set-TextValue $ws 'D11' '0.07951'
Set-TextValue $ws 'D12' '97.39'
Set-TextValue $ws 'D13' '1.859.11'
Set-TextValue $ws 'D14' '5.081'
Set-TextValue $ws 'D15' '0.6739'
Set-TextValue $ws 'D16' '265.74'
Set-TextValue $ws 'D17' '30.151.08'
Set-TextValue $ws 'D19' '1.002'
Set-TextValue $ws 'D20' '0.000007539'
Set-TextValue $ws 'D21' '2.101.23'
Set-TextValue $ws 'D23' '5.215'
Set-TextValue $ws 'D24' '6.133'
Set-TextValue $ws 'D25' '166.53'
Set-TextValue $ws 'D26' '9.145'
Set-TextValue $ws 'D28' '1.924'
Set-TextValue $ws 'D29' '1.395'
Set-TextValue $ws 'D30' '0.09825'
Set-TextValue $ws 'D31' '1.465'
Set-TextValue $ws 'D32' '4.261'
Set-TextValue $ws 'D33' '3.988'
Set-TextValue $ws 'D34' '0.04683'
Set-TextValue $ws 'D36' '0.6963'
Set-TextValue $ws 'D38' '0.01858'
Set-TextValue $ws 'D39' '2.602'
Set-TextValue $ws 'D40' '6.315'
Set-TextValue $ws 'D41' '73.23'
Set-TextValue $ws 'D42' '1.928'
Set-TextValue $ws 'D43' '1.002'
Set-TextValue $ws 'D44' '0.8358'
Set-TextValue $ws 'D45' '103.10'
Set-TextValue $ws 'D46' '0.4116'
Set-TextValue $ws 'D47' '941.58'
Set-TextValue $ws 'D48' '9.094'
Set-TextValue $ws 'D49' '6.994'
Set-TextValue $ws 'D50' '33.74'

# Column E (Volume/1h) updates
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('E8').Value = '  +1.83%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('E10').Value = '  +1.54%  '
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('E16').Value = '  -5.20%  '
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('E18').Value = '  +7.27%  '
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('E20').Value = '  +3.77%  '
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('E23').Value = '  -4.78%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('E29').Value = '  +1.45%  '
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('E32').Value = '  -2.96%  '
$ws.Range('E33').Value = '  -2.63%  '
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('E39').Value = '  +2.78%  '
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('E46').Value = '  -1.35%  '
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('E49').Value = '  -2.61%  '
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('E51').Value = '  +0.43%  '
